$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.625.65"
$ws.Range("E2").Value = "  +2.14%  "
$ws.Range("D3").Value = "1.888.87"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'245.75"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "'0.4910"
$ws.Range("E7").Value = "  -0.55%  "
$ws.Range("D8").Value = "'0.2945"
$ws.Range("D9").Value = "'0.06763"
$ws.Range("E9").Value = "  +1.81%  "
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("D11").Value = "'17.18"
$ws.Range("E11").Value = "  +2.47%  "
$ws.Range("D12").Value = "'0.07245"
$ws.Range("E12").Value = "  +0.90%  "
$ws.Range("D13").Value = "'91.03"
$ws.Range("E13").Value = "  +4.92%  "
$ws.Range("D14").Value = "'0.6769"
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("D15").Value = "'5.050"
$ws.Range("E15").Value = "  +3.18%  "
$ws.Range("D16").Value = "30.614.42"
$ws.Range("E16").Value = "  +2.18%  "
$ws.Range("D17").Value = "'0.000007950"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").Value = "'13.14"
$ws.Range("E19").Value = "  +2.63%  "
$ws.Range("D20").Value = "2.130.55"
$ws.Range("E20").Value = "  +0.39%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "'4.821"
$ws.Range("E22").Value = "  +0.68%  "
$ws.Range("D23").Value = "'185.19"
$ws.Range("E23").Value = "  +30.22%  "
$ws.Range("D24").Value = "'6.074"
$ws.Range("E24").Value = "  +3.10%  "
$ws.Range("D25").Value = "'9.336"
$ws.Range("E25").Value = "  +2.28%  "
$ws.Range("E26").Value = "  +3.18%  "
$ws.Range("D27").Value = "'19.06"
$ws.Range("E27").Value = "  +11.71%  "
$ws.Range("D28").Value = "'1.901"
$ws.Range("E28").Value = "  -1.12%  "
$ws.Range("D29").Value = "'1.397"
$ws.Range("E29").Value = "  +0.70%  "
$ws.Range("D30").Value = "'4.340"
$ws.Range("E30").Value = "  +2.83%  "
$ws.Range("E31").Value = "  +3.02%  "
$ws.Range("D32").Value = "'3.998"
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("D33").Value = "'0.05187"
$ws.Range("E33").Value = "  +2.59%  "
$ws.Range("D34").Value = "'0.7534"
$ws.Range("E34").Value = "  +4.90%  "
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("D36").Value = "'2.750"
$ws.Range("E36").Value = "  +3.17%  "
$ws.Range("E37").Value = "  +3.06%  "
$ws.Range("D38").Value = "'2.661"
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("D39").Value = "'2.140"
$ws.Range("E39").Value = "  -1.61%  "
$ws.Range("D40").Value = "'0.9357"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").Value = "'0.4412"
$ws.Range("E41").Value = "  +3.91%  "
$ws.Range("D42").Value = "'105.28"
$ws.Range("E42").Value = "  +1.23%  "
$ws.Range("D43").Value = "'1.000"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").Value = "'5.748"
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").Value = "'7.577"
$ws.Range("E45").Value = "  +2.31%  "
$ws.Range("E46").Value = "  +4.86%  "
$ws.Range("D47").Value = "'0.05848"
$ws.Range("E47").Value = "  +2.59%  "
$ws.Range("E48").Value = "  +7.55%  "
$ws.Range("D49").Value = "'8.695"
$ws.Range("E49").Value = "  +4.90%  "
$ws.Range("D50").Value = "'0.3927"
$ws.Range("E50").Value = "  +3.82%  "
$ws.Range("D51").Value = "'33.49"
$ws.Range("E51").Value = "  +2.54%  "

# Reset style on cells that were text-forced via quote-prefix so no stray
# quotePrefix style attribute is left on the cell (matches original plain style).
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
